# Update column G ("K") values for rows 2-18 per the diff.
# Old Strike# values are replaced with newly computed K values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 5
    3  = 2
    4  = 3
    5  = 6
    6  = 5
    7  = 6
    8  = 2
    9  = 6
    10 = 7
    11 = 2
    12 = 2
    13 = 5
    14 = 5
    15 = 0
    16 = 1
    17 = 0
    18 = 1
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
